$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Sending=ECs, Target=ECs
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Ncam1"
$ws.Cells.Item(2,3).Value = "Robo1"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 0.7004376666666666
$ws.Cells.Item(2,8).Value = 2.101313
$ws.Cells.Item(2,9).Value = 0.04511966030063898
$ws.Cells.Item(2,10).Value = 0.04511966030063898
$ws.Cells.Item(2,11).Value = 2
$ws.Cells.Item(2,12).Value = 0.6666666666666666
$ws.Cells.Item(2,13).Value = 0.1757713333333334
$ws.Cells.Item(2,14).Value = 0.5273140000000001
$ws.Cells.Item(2,15).Value = 0.009444264308298454
$ws.Cells.Item(2,16).Value = 0.009444264308298454
$ws.Cells.Item(2,17).Value = 0.1231168625868889
$ws.Cells.Item(2,18).Value = 1.108051763282
$ws.Cells.Item(2,19).Value = 0.0004261219973798754
$ws.Cells.Item(2,20).Value = 0.0004261219973798754

# Row 3: Sending=ECs, Target=FAPs
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Ncam1"
$ws.Cells.Item(3,3).Value = "Robo1"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 0.7004376666666666
$ws.Cells.Item(3,8).Value = 2.101313
$ws.Cells.Item(3,9).Value = 0.04511966030063898
$ws.Cells.Item(3,10).Value = 0.04511966030063898
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 16.98312366666667
$ws.Cells.Item(3,14).Value = 50.949371
$ws.Cells.Item(3,15).Value = 0.9125100529581165
$ws.Cells.Item(3,16).Value = 0.9125100529581165
$ws.Cells.Item(3,17).Value = 11.89561951379144
$ws.Cells.Item(3,18).Value = 107.060575624123
$ws.Cells.Item(3,19).Value = 0.0411721436103883
$ws.Cells.Item(3,20).Value = 0.0411721436103883

# Row 4: Sending=ECs, Target=sCs
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Ncam1"
$ws.Cells.Item(4,3).Value = "Robo1"
$ws.Cells.Item(4,4).Value = "sCs"
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 0.7004376666666666
$ws.Cells.Item(4,8).Value = 2.101313
$ws.Cells.Item(4,9).Value = 0.04511966030063898
$ws.Cells.Item(4,10).Value = 0.04511966030063898
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 1.452542333333333
$ws.Cells.Item(4,14).Value = 4.357627
$ws.Cells.Item(4,15).Value = 0.07804568273358503
$ws.Cells.Item(4,16).Value = 0.07804568273358505
$ws.Cells.Item(4,17).Value = 1.017415362694555
$ws.Cells.Item(4,18).Value = 9.156738264250999
$ws.Cells.Item(4,19).Value = 0.003521394692870802
$ws.Cells.Item(4,20).Value = 0.003521394692870802

# Row 5: Sending=FAPs, Target=ECs
$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "Ncam1"
$ws.Cells.Item(5,3).Value = "Robo1"
$ws.Cells.Item(5,4).Value = "ECs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 0.237305
$ws.Cells.Item(5,8).Value = 0.711915
$ws.Cells.Item(5,9).Value = 0.01528632952964618
$ws.Cells.Item(5,10).Value = 0.01528632952964618
$ws.Cells.Item(5,11).Value = 2
$ws.Cells.Item(5,12).Value = 0.6666666666666666
$ws.Cells.Item(5,13).Value = 0.1757713333333334
$ws.Cells.Item(5,14).Value = 0.5273140000000001
$ws.Cells.Item(5,15).Value = 0.009444264308298454
$ws.Cells.Item(5,16).Value = 0.009444264308298454
$ws.Cells.Item(5,17).Value = 0.04171141625666667
$ws.Cells.Item(5,18).Value = 0.37540274631
$ws.Cells.Item(5,19).Value = 0.0001443681363817261
$ws.Cells.Item(5,20).Value = 0.0001443681363817261

# Row 6: Sending=FAPs, Target=FAPs
$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Ncam1"
$ws.Cells.Item(6,3).Value = "Robo1"
$ws.Cells.Item(6,4).Value = "FAPs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 0.237305
$ws.Cells.Item(6,8).Value = 0.711915
$ws.Cells.Item(6,9).Value = 0.01528632952964618
$ws.Cells.Item(6,10).Value = 0.01528632952964618
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 16.98312366666667
$ws.Cells.Item(6,14).Value = 50.949371
$ws.Cells.Item(6,15).Value = 0.9125100529581165
$ws.Cells.Item(6,16).Value = 0.9125100529581165
$ws.Cells.Item(6,17).Value = 4.030180161718333
$ws.Cells.Item(6,18).Value = 36.271621455465
$ws.Cells.Item(6,19).Value = 0.01394892936863265
$ws.Cells.Item(6,20).Value = 0.01394892936863265

# Row 7: Sending=FAPs, Target=sCs
$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Ncam1"
$ws.Cells.Item(7,3).Value = "Robo1"
$ws.Cells.Item(7,4).Value = "sCs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 0.237305
$ws.Cells.Item(7,8).Value = 0.711915
$ws.Cells.Item(7,9).Value = 0.01528632952964618
$ws.Cells.Item(7,10).Value = 0.01528632952964618
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 1.452542333333333
$ws.Cells.Item(7,14).Value = 4.357627
$ws.Cells.Item(7,15).Value = 0.07804568273358503
$ws.Cells.Item(7,16).Value = 0.07804568273358505
$ws.Cells.Item(7,17).Value = 0.3446955584116667
$ws.Cells.Item(7,18).Value = 3.102260025705
$ws.Cells.Item(7,19).Value = 0.001193032024631798
$ws.Cells.Item(7,20).Value = 0.001193032024631798

# Row 8: Sending=sCs, Target=ECs
$ws.Cells.Item(8,1).Value = "sCs"
$ws.Cells.Item(8,2).Value = "Ncam1"
$ws.Cells.Item(8,3).Value = "Robo1"
$ws.Cells.Item(8,4).Value = "ECs"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 14.58625866666667
$ws.Cells.Item(8,8).Value = 43.758776
$ws.Cells.Item(8,9).Value = 0.9395940101697148
$ws.Cells.Item(8,10).Value = 0.9395940101697148
$ws.Cells.Item(8,11).Value = 2
$ws.Cells.Item(8,12).Value = 0.6666666666666666
$ws.Cells.Item(8,13).Value = 0.1757713333333334
$ws.Cells.Item(8,14).Value = 0.5273140000000001
$ws.Cells.Item(8,15).Value = 0.009444264308298454
$ws.Cells.Item(8,16).Value = 0.009444264308298454
$ws.Cells.Item(8,17).Value = 2.563846134184889
$ws.Cells.Item(8,18).Value = 23.074615207664
$ws.Cells.Item(8,19).Value = 0.008873774174536851
$ws.Cells.Item(8,20).Value = 0.008873774174536851

# Row 9: Sending=sCs, Target=FAPs
$ws.Cells.Item(9,1).Value = "sCs"
$ws.Cells.Item(9,2).Value = "Ncam1"
$ws.Cells.Item(9,3).Value = "Robo1"
$ws.Cells.Item(9,4).Value = "FAPs"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 14.58625866666667
$ws.Cells.Item(9,8).Value = 43.758776
$ws.Cells.Item(9,9).Value = 0.9395940101697148
$ws.Cells.Item(9,10).Value = 0.9395940101697148
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 16.98312366666667
$ws.Cells.Item(9,14).Value = 50.949371
$ws.Cells.Item(9,15).Value = 0.9125100529581165
$ws.Cells.Item(9,16).Value = 0.9125100529581165
$ws.Cells.Item(9,17).Value = 247.7202347699885
$ws.Cells.Item(9,18).Value = 2229.482112929896
$ws.Cells.Item(9,19).Value = 0.8573889799790955
$ws.Cells.Item(9,20).Value = 0.8573889799790955

# Row 10: Sending=sCs, Target=sCs
$ws.Cells.Item(10,1).Value = "sCs"
$ws.Cells.Item(10,2).Value = "Ncam1"
$ws.Cells.Item(10,3).Value = "Robo1"
$ws.Cells.Item(10,4).Value = "sCs"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 14.58625866666667
$ws.Cells.Item(10,8).Value = 43.758776
$ws.Cells.Item(10,9).Value = 0.9395940101697148
$ws.Cells.Item(10,10).Value = 0.9395940101697148
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 1.452542333333333
$ws.Cells.Item(10,14).Value = 4.357627
$ws.Cells.Item(10,15).Value = 0.07804568273358503
$ws.Cells.Item(10,16).Value = 0.07804568273358505
$ws.Cells.Item(10,17).Value = 21.18715819828356
$ws.Cells.Item(10,18).Value = 190.684423784552
$ws.Cells.Item(10,19).Value = 0.07333125601608242
$ws.Cells.Item(10,20).Value = 0.07333125601608244
